$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 769, shifting existing rows 769-818 down to 770-819
$ws.Rows.Item(769).Insert()

# Populate the newly inserted row 769 with its data
$ws.Cells.Item(769, 1).Value = 5
$ws.Cells.Item(769, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(769, 3).Value = "Maule"
$ws.Cells.Item(769, 4).Value = 45021
$ws.Cells.Item(769, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(769, 5).Value = 7
$ws.Cells.Item(769, 6).Value = 100112002
$ws.Cells.Item(769, 7).Value = "Pimiento"
$ws.Cells.Item(769, 8).Value = "Zafiro rojo"
$ws.Cells.Item(769, 9).Value = "Primera"
$ws.Cells.Item(769, 10).Value = 200
$ws.Cells.Item(769, 11).Value = 15000
$ws.Cells.Item(769, 12).Value = 15000
$ws.Cells.Item(769, 13).Value = 15000
$ws.Cells.Item(769, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(769, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(769, 16).Value = 1000
$ws.Cells.Item(769, 17).Value = 15
$ws.Cells.Item(769, 18).Value = "Hortaliza"
